$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at position 85, pushing the existing
# rows 85-123 down to 86-124 (dimension grows from A1:T123 to A1:T124).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new record
# (Vega Modelo de Temuco / Durazno / Early Majestic).
$ws.Range('A85').Value = 10
$ws.Range('B85').Value = 'Vega Modelo de Temuco'
$ws.Range('C85').Value = 'La Araucanía'
$ws.Range('D85').Value = 44488
$ws.Range('E85').Value = 9
$ws.Range('F85').Value = 'Fruta'
$ws.Range('G85').Value = 100103
$ws.Range('H85').Value = 'Frutos de hueso (carozo)'
$ws.Range('I85').Value = 100103004
$ws.Range('J85').Value = 'Durazno'
$ws.Range('K85').Value = 'Early Majestic'
$ws.Range('L85').Value = 'Primera'
$ws.Range('M85').Value = 80
$ws.Range('N85').Value = 20000
$ws.Range('O85').Value = 20000
$ws.Range('P85').Value = 20000
$ws.Range('Q85').Value = '$/bandeja 10 kilos granel'
$ws.Range('R85').Value = 'Provincia de Limarí'
$ws.Range('S85').Value = 2000
$ws.Range('T85').Value = 10
